# Natmi following Dr Hou advice
# Update row 2 values and add a new row 3 (sCs -> Ccl21b/Ackr2 -> FAPs) pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing FAPs -> Ccl21b/Ackr2 -> FAPs edge with recalculated stats ---
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.380435
$ws.Range("H2").Value = 1.141305
$ws.Range("I2").Value = 0.7997108917301441
$ws.Range("J2").Value = 0.7997108917301442
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.35607266666667
$ws.Range("N2").Value = 37.068218
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 4.700682504943334
$ws.Range("R2").Value = 42.30614254449
$ws.Range("S2").Value = 0.7997108917301441
$ws.Range("T2").Value = 0.7997108917301442

# --- Row 3: new sCs -> Ccl21b/Ackr2 -> FAPs edge ---
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.09528066666666667
$ws.Range("H3").Value = 0.285842
$ws.Range("I3").Value = 0.2002891082698559
$ws.Range("J3").Value = 0.2002891082698559
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.35607266666667
$ws.Range("N3").Value = 37.068218
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.177294841061778
$ws.Range("R3").Value = 10.595653569556
$ws.Range("S3").Value = 0.2002891082698559
$ws.Range("T3").Value = 0.2002891082698559
